$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "2025/11/26"
$ws.Range("I1").Style = $ws.Range("H1").Style

$ws.Range("I2").Value = 59
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = 59
$ws.Range("I5").Value = 59
$ws.Range("I6").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("I8").Value = 59
$ws.Range("I9").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("I11").Value = 59
$ws.Range("I12").Value = ""
$ws.Range("I13").Value = ""
$ws.Range("I14").Value = ""
$ws.Range("I15").Value = 59
$ws.Range("I16").Value = ""
$ws.Range("I17").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("I19").Value = ""
$ws.Range("I20").Value = 0.6
$ws.Range("I21").Value = ""
$ws.Range("I22").Value = ""
$ws.Range("I23").Value = 0.6
$ws.Range("I24").Value = 0.59
$ws.Range("I25").Value = 0.59
